# The diff inserts one new data row (a new weekly "Perejil" price record for
# Feria Lagunitas de Puerto Montt) immediately above what was previously row 52.
# Every subsequent row (old 52..168) shifts down by one (to 53..169), and the
# sheet's used range grows from A1:R168 to A1:R169.
#
# Insert a whole row at row 52 to push existing data down, then populate the
# freshly inserted row with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(52).Insert()

$ws.Range("A52").Value = 4
$ws.Range("B52").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C52").Value = "Los Lagos"
$ws.Range("D52").Value = 44498
$ws.Range("E52").Value = 10
$ws.Range("F52").Value = 100112044
$ws.Range("G52").Value = "Perejil"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 160
$ws.Range("K52").Value = 4500
$ws.Range("L52").Value = 4500
$ws.Range("M52").Value = 4500
$ws.Range("N52").Value = "$/docena de atados (3 kilos)"
$ws.Range("O52").Value = "Región Metropolitana"
$ws.Range("P52").Value = 1500
$ws.Range("Q52").Value = 3
$ws.Range("R52").Value = "Hortaliza"
